$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows after the last existing data row (28). Inserting with a
# shift-down picks up the formatting of the row above, which reproduces the
# per-column cell styles (date / text / text / text / number) used throughout
# the table.
$ws.Rows("29:30").Insert(-4121) | Out-Null   # xlShiftDown

$dataDateSerial = 43242   # 22/05/2018

# Row 29: Licata Rosa / Spilii 250 gr. / N°. / 1
$ws.Range("A29").Value = $dataDateSerial
$ws.Range("B29").Value = "Licata Rosa"
$ws.Range("C29").Value = "Spilii 250 gr."
$ws.Range("D29").Value = "N°."
$ws.Range("E29").Value = 1

# Row 30: Licata Rosa / Centimetro / N°. / 2
$ws.Range("A30").Value = $dataDateSerial
$ws.Range("B30").Value = "Licata Rosa"
$ws.Range("C30").Value = "Centimetro"
$ws.Range("D30").Value = "N°."
$ws.Range("E30").Value = 2
